$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the "Currency" label / value row (row 6) on the input sheet:
#  - label text is lower-cased to match the other field-name cells
#  - value text drops its trailing space
#  - the value cell picks up a plain green-fill style
$wsInput.Range("A6").Value = "currency"
$wsInput.Range("B6").Value = "US Dollar"
$wsInput.Range("B6").Style = "Normal"
$wsInput.Range("B6").Interior.Color = 5296274

# Move the active tab / selection back to the input sheet
$wsInput.Activate()
$wsInput.Range("B10").Select()
